# Add the GitHub project-link hyperlink text to the "Content Placeholder 2"
# shape on the "PROJECT LINK" slide (sldId 264 / cId 1430217275, shape id 3,
# creationId {54BCE890-181B-973C-F6CA-E2F9969E94E3}).
#
# The target text "dsk2711/keylooger (github.com)" is inserted as three runs
# (matching how PowerPoint splits runs when text is typed/auto-corrected
# incrementally), each one hyperlinked to the GitHub repo URL.

$p = $ppt.ActivePresentation

# Locate the slide by its SlideID (264) rather than assuming a fixed index.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 264) {
        $targetSlide = $p.Slides.Item($i)
        break
    }
}

# Locate the shape by its Id (3) within that slide.
$targetShape = $null
for ($j = 1; $j -le $targetSlide.Shapes.Count; $j++) {
    if ($targetSlide.Shapes.Item($j).Id -eq 3) {
        $targetShape = $targetSlide.Shapes.Item($j)
        break
    }
}

$tr = $targetShape.TextFrame.TextRange
$url = "https://github.com/dsk2711/keylooger"

$r1 = $tr.InsertAfter("dsk2711/")
$r1.LanguageID = "en-IN"
$r1.ActionSettings.Item(1).Hyperlink.Address = $url

$r2 = $tr.InsertAfter("keylooger")
$r2.LanguageID = "en-IN"
$r2.ActionSettings.Item(1).Hyperlink.Address = $url

$r3 = $tr.InsertAfter(" (github.com)")
$r3.LanguageID = "en-IN"
$r3.ActionSettings.Item(1).Hyperlink.Address = $url

Write-Host "Updated shape text: [$($tr.Text)]"
